$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9725.4
$ws.Range("I69").Value = 7988.5
$ws.Range("K69").Value = 23965.5
$ws.Range("M69").Value = -23091.5
$ws.Range("H72").Value = 9725.4
$ws.Range("I72").Value = 7988.5
$ws.Range("K72").Value = 71896.5
$ws.Range("M72").Value = -67528.5
$ws.Range("H98").Value = 2208.818
$ws.Range("I98").Value = 2208.818
$ws.Range("K98").Value = 2208.818
$ws.Range("M98").Value = -710.8180000000002
$ws.Range("H113").Value = 6583
$ws.Range("I113").Value = 2750
$ws.Range("K113").Value = 2750
$ws.Range("M113").Value = 504
$ws.Range("H116").Value = 34039
$ws.Range("I116").Value = 52998.332
$ws.Range("K116").Value = 52998.332
$ws.Range("M116").Value = -49556.332
$ws.Range("H122").Value = 2208.818
$ws.Range("I122").Value = 2208.818
$ws.Range("K122").Value = 6626.454000000001
$ws.Range("M122").Value = -4176.454000000001
$ws.Range("H138").Value = 3461.5354
$ws.Range("I138").Value = 1749.5454
$ws.Range("J138").Value = 3950.6753
$ws.Range("K138").Value = 5248.6362
$ws.Range("L138").Value = 11852.0259
$ws.Range("M138").Value = -108.6361999999999
$ws.Range("N138").Value = -22132.0259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3161.875
$ws.Range("I45").Value = 2508.182
$ws.Range("J45").Value = 4600
$ws.Range("K45").Value = 2508.182
$ws.Range("L45").Value = 4600
$ws.Range("M45").Value = -2131.182
$ws.Range("N45").Value = -5354
$ws.Range("H74").Value = 3609
$ws.Range("I74").Value = 3696
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 3696
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -2822
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 3609
$ws.Range("I77").Value = 3696
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 18480
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -14112
$ws.Range("N77").Value = -23736
$ws.Range("H88").Value = 2999.5
$ws.Range("J88").Value = 2999.5
$ws.Range("L88").Value = 2999.5
$ws.Range("N88").Value = -3811.5
$ws.Range("H91").Value = 2999.5
$ws.Range("J91").Value = 2999.5
$ws.Range("L91").Value = 2999.5
$ws.Range("N91").Value = -5807.5
$ws.Range("H102").Value = 2779.5386
$ws.Range("I102").Value = 2683.4
$ws.Range("J102").Value = 3100
$ws.Range("K102").Value = 2683.4
$ws.Range("L102").Value = 3100
$ws.Range("M102").Value = -1061.4
$ws.Range("N102").Value = -6344

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 245.8
$ws.Range("I5").Value = 177
$ws.Range("J5").Value = 291.66666
$ws.Range("K5").Value = 177
$ws.Range("L5").Value = 291.66666
$ws.Range("M5").Value = -64
$ws.Range("N5").Value = -517.66666
$ws.Range("H105").Value = 2696.1765
$ws.Range("I105").Value = 2389
$ws.Range("K105").Value = 2389
$ws.Range("M105").Value = -642

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1809.069
$ws.Range("I31").Value = 1267.4333
$ws.Range("K31").Value = 1267.4333
$ws.Range("M31").Value = -972.4332999999999
$ws.Range("H34").Value = 1809.069
$ws.Range("I34").Value = 1267.4333
$ws.Range("K34").Value = 1267.4333
$ws.Range("M34").Value = -1065.4333
$ws.Range("H41").Value = 64999.75
$ws.Range("I41").Value = 45000
$ws.Range("K41").Value = 45000
$ws.Range("M41").Value = -44572
$ws.Range("H53").Value = 29995
$ws.Range("J53").Value = 29995
$ws.Range("L53").Value = 29995
$ws.Range("N53").Value = -31209
$ws.Range("H99").Value = 2929.9285
$ws.Range("I99").Value = 2551.9
$ws.Range("K99").Value = 2551.9
$ws.Range("M99").Value = -1053.9
$ws.Range("H104").Value = 79890
$ws.Range("J104").Value = 79890
$ws.Range("L104").Value = 79890
$ws.Range("N104").Value = -85132
$ws.Range("H126").Value = 2929.9285
$ws.Range("I126").Value = 2551.9
$ws.Range("K126").Value = 7655.700000000001
$ws.Range("M126").Value = -5185.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 62.583332
$ws.Range("I12").Value = 9.666667
$ws.Range("K12").Value = 29.000001
$ws.Range("M12").Value = 143.999999
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H87").Value = 6969
$ws.Range("I87").Value = 6969
$ws.Range("K87").Value = 20907
$ws.Range("M87").Value = -19659
$ws.Range("H90").Value = 6969
$ws.Range("I90").Value = 6969
$ws.Range("K90").Value = 62721
$ws.Range("M90").Value = -56481
$ws.Range("H131").Value = 2660.6428
$ws.Range("I131").Value = 9666
$ws.Range("J131").Value = 2121.7693
$ws.Range("K131").Value = 28998
$ws.Range("L131").Value = 6365.3079
$ws.Range("M131").Value = -23958
$ws.Range("N131").Value = -16445.3079
$ws.Range("H140").Value = 980.6667
$ws.Range("I140").Value = 721
$ws.Range("K140").Value = 2163
$ws.Range("M140").Value = 3017

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1996.8823
$ws.Range("I22").Value = 1642.6428
$ws.Range("J22").Value = 3650
$ws.Range("K22").Value = 1642.6428
$ws.Range("L22").Value = 3650
$ws.Range("M22").Value = -1347.6428
$ws.Range("N22").Value = -4240
$ws.Range("H27").Value = 1996.8823
$ws.Range("I27").Value = 1642.6428
$ws.Range("J27").Value = 3650
$ws.Range("K27").Value = 1642.6428
$ws.Range("L27").Value = 3650
$ws.Range("M27").Value = -1535.6428
$ws.Range("N27").Value = -3864
$ws.Range("H46").Value = 2931.8
$ws.Range("I46").Value = 1206.5
$ws.Range("J46").Value = 3260.4285
$ws.Range("K46").Value = 1206.5
$ws.Range("L46").Value = 3260.4285
$ws.Range("M46").Value = -1018.5
$ws.Range("N46").Value = -3636.4285
$ws.Range("H55").Value = 1776
$ws.Range("J55").Value = 1114.3334
$ws.Range("L55").Value = 1114.3334
$ws.Range("N55").Value = -1460.3334
$ws.Range("H93").Value = 66668028
$ws.Range("I93").Value = 76924190
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 76924190
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -76922942
$ws.Range("N93").Value = -5496
$ws.Range("H111").Value = 63747.5
$ws.Range("J111").Value = 63747.5
$ws.Range("L111").Value = 63747.5
$ws.Range("N111").Value = -71927.5
$ws.Range("H132").Value = 19664.834
$ws.Range("I132").Value = 19997.8
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 59993.39999999999
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -57463.39999999999
$ws.Range("N132").Value = -59060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4471.174
$ws.Range("I81").Value = 972.1539
$ws.Range("J81").Value = 9019.9
$ws.Range("K81").Value = 1944.3078
$ws.Range("L81").Value = 18039.8
$ws.Range("M81").Value = -883.3078
$ws.Range("N81").Value = -20161.8
$ws.Range("H84").Value = 4471.174
$ws.Range("I84").Value = 972.1539
$ws.Range("J84").Value = 9019.9
$ws.Range("K84").Value = 9721.539000000001
$ws.Range("L84").Value = 90199
$ws.Range("M84").Value = -4417.539000000001
$ws.Range("N84").Value = -100807
$ws.Range("H115").Value = 101965.4
$ws.Range("J115").Value = 101965.4
$ws.Range("L115").Value = 101965.4
$ws.Range("N115").Value = -105099.4
$ws.Range("H129").Value = 125998
$ws.Range("J129").Value = 125998
$ws.Range("L129").Value = 125998
$ws.Range("N129").Value = -135998
$ws.Range("H132").Value = 1683.1578
$ws.Range("I132").Value = 1795.9697
$ws.Range("K132").Value = 5387.909100000001
$ws.Range("M132").Value = -2857.909100000001
